# north-virginia-2024.xlsx edits:
#  - a few county "result" values changed
#  - the view was scrolled down and a new cell selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- updated county values (column C) ---
$ws.Range("C18").Value = -53.5
$ws.Range("C31").Value = -64.1
$ws.Range("C50").Value = -57.3
$ws.Range("C65").Value = 56.7

# --- updated sheet view: scrolled so row 46 is at the top, C66 selected ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C66").Select()
